$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "spectraltype(short)" column (column C) values: the abbreviation
# now includes the luminosity class (e.g. "B1" -> "B1I"), and rows whose
# spectral type is a giant (class III) show the full short spectral type
# instead of just adding an "I" suffix.
$ws.Range("C2").Value  = "B1I"
$ws.Range("C3").Value  = "O9I"
$ws.Range("C4").Value  = "O8III"
$ws.Range("C5").Value  = "B0I"
$ws.Range("C6").Value  = "O9III"
$ws.Range("C7").Value  = "B1I"
$ws.Range("C8").Value  = "B0I"
$ws.Range("C9").Value  = "O6I"
$ws.Range("C10").Value = "O8I"
$ws.Range("C11").Value = "O8I"
$ws.Range("C12").Value = "O9I"

# Keep the same selected cell as the saved workbook (moved one row down).
[void]$ws.Range("C14").Select()
